$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 38 (Idaho) turned into an error row for this run:
# - Data columns B..H get cleared (no value, no date style)
# - J38 boolean flips from TRUE to FALSE
# - O38 status message becomes the timeout error text

$ws.Range("B38:H38").Value = ""
$ws.Range("B38").ClearFormats()

$ws.Range("J38").Value = $false

$ws.Range("O38").Value = "An error occurred. ... TimeoutException('', None, None)"
